# Reorders the "Recorded By" (column G) comma-separated list so that the
# literal token "System" (capital S) is moved from its original position to
# the end of the list, for every row where it appears. All other tokens keep
# their original relative order. Rows whose G value does not contain the
# exact token "System" are left untouched.

function Reorder-RecordedBy {
    param([string]$val)

    $parts = $val -split ", "

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }
    if (-not $hasSystem) { return $val }

    $result = ""
    $first = $true
    foreach ($p in $parts) {
        if (-not $p.Equals("System")) {
            if ($first) {
                $result = $p
                $first = $false
            } else {
                $result = $result + ", " + $p
            }
        }
    }
    if ($first) {
        $result = "System"
    } else {
        $result = $result + ", System"
    }
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text
    if ($current -ne $null -and $current -ne "") {
        $updated = Reorder-RecordedBy $current
        if (-not $updated.Equals($current)) {
            $cell.Value = $updated
        }
    }
}
